# Binary growth results - add CGXII ("CGXlab+") growth data sheet
$wb = $excel.ActiveWorkbook

# Remember the current active sheet (CGXlab) selection before we touch anything,
# so the recorded cursor position for that tab survives it losing focus.
$ws4 = $wb.Worksheets.Item("CGXlab")
$ws4.Range("F18").Select()

# Add the new worksheet at the very end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add($null, $lastSheet)
$new.Name = "CGXlab+"

# Column width for column A (strain / medium names).
$new.Columns.Item(1).ColumnWidth = 31.83

# Fill cell-by-cell in the same order the original data was typed in
# (this keeps the shared-string table's insertion order faithful).

# Column A: strain / medium names, top to bottom.
$new.Range("A1").Value = "medium"
$new.Range("A2").Value = "CGXII"
$new.Range("A3").Value = "CGXII+cobalt,pnto,cys__L"
$new.Range("A4").Value = "CGXII-ni,PCA,btn"
$new.Range("A5").Value = "CGXII-ni,PCA,btn+cobalt,pnto,cys__L"

# Column B: modus tags.
$new.Range("B2").Value = "complete"
$new.Range("B5").Value = "minimal"
$new.Range("B4").Value = "deprived"
$new.Range("B3").Value = "full"

# Column C/D/E numeric data plus the header row labels.
$new.Range("C1").Value = "16-1"
$new.Range("D1").Value = "16-2"
$new.Range("E1").Value = "16-2"
$new.Range("F1").Value = 16

$new.Range("C2").Value = 0.155
$new.Range("D2").Value = 0.166
$new.Range("E2").Value = 0.219

$new.Range("C3").Value = 0.224
$new.Range("D3").Value = 0.263
$new.Range("E3").Value = 0.253

$new.Range("C4").Value = 0.127
$new.Range("D4").Value = 0.149
$new.Range("E4").Value = 0.125

$new.Range("C5").Value = 0.185
$new.Range("D5").Value = 0.214
$new.Range("E5").Value = 0.207

# Columns H/I/J: additive supplement reference table.
$new.Range("H1").Value = "cobalt"
$new.Range("H2").Value = "0.013 g/L"
$new.Range("I1").Value = "pnto__R"
$new.Range("I2").Value = "0.001 g/L"
$new.Range("J1").Value = "cys__L"
$new.Range("J2").Value = "0.121 g/L"

# ---- Column F: average growth. F2 stands alone; F3:F5 share one formula
#      (mirrors the E/F columns on the CGXlab sheet, where row 2 is also
#      the odd one out).
$new.Range("F2").Formula = "=SUM(C2:E2)/3"
$new.Range("F3:F5").Formula = "=SUM(C3:E3)/3"
$new.Range("F2:F5").NumberFormat = "0.000"

# Park the cursor where the author left it on the new sheet.
$new.Range("C26").Select()
